$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values (fixes to existing product rows)
$ws.Range("D2").Value = 189
$ws.Range("C16").Value = 148000

$ws.Range("C40").Value = 25000
$ws.Range("D40").Value = 1
$ws.Range("J40").Value = 2

$ws.Range("C41").Value = 70000
$ws.Range("D41").Value = 1
$ws.Range("J41").Value = 2

# Remove the last two product rows (42 and 43) entirely
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(42).Delete()
